$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Compressed Size" to "C.Size"
$ws.Name = "C.Size"

# Update the title text in A1 (shorten it, dropping the C.Time/Avg. C. Ratio mention)
$ws.Range("A1").Value = "The C.Size (B) of the compressed file on two highly repetitive RNA sequence obtained by GraSS and other benchmark methods"

# Update the active selection shown when the workbook is opened
$ws.Range("D16").Select()
